$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (shifts N->O, O->P, P->Q),
# matching the structural column-insert reflected in the diff.
[void]$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab and update its selection,
# which also clears the previous tabSelected flag on "Transactions".
$ws.Activate()
[void]$ws.Range("R6").Select()
